$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Validation")

$data = @(
    @("Table","Field","Type","Read/Write","Data Source","Constraints","TestID","Pass/Fail"),
    @("User","Name","varchar ","write","sForm.textInput1","String <255","data.1","Pass"),
    @("User","Email","varchar ","write","sForm.textInput1","String <255","data.2","Pass"),
    @("User","Password","varchar ","write","sForm.textInput1","String > 9","data.3","Pass"),
    @("User","RecipeList","ManytoMany Filed (Django)","write","Sform.writeSQL","None","data.4","Pass"),
    @("User's Saved Recipes","User(Email)","varchar ","read","Sform.writeSQL","String <255","data.5","Fail"),
    @("User's Saved Recipes","Recipe_id","int","read","Sform.writeSQL","None","data.6","Fail"),
    @("User's Saved Recipes","relation_id","int","write","Sform.writeSQL","None","data.7","Fail"),
    @("Recipes","id","int","read","Sform.writeSQL","None","data.8","Pass"),
    @("Recipes","ingredient","varchar ","write","sForm.textInput1","None","data.9","Pass"),
    @("Recipes","Recipe","varchar ","read","Sform.writeSQL","None","data.10","Pass"),
    @("Recipes","image","varchar (image_field)","read","Sform.writeSQL","None","data.11","Pass")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws.Range("I5").Select()

$ws.Columns.Item(1).ColumnWidth = 17.42
$ws.Columns.Item(3).ColumnWidth = 22.75
$ws.Columns.Item(5).ColumnWidth = 13.75

$wsRemediation = $wb.Worksheets.Item("Remediation")
$wsRemediation.Activate()
